# TC_003_Registration Module - UK Account Registration test scripts
# Adds UK_LoginDatas and USA_LoginDatas sheets (mirroring India_LoginDatas),
# refreshes India_LoginDatas' own test rows/hyperlinks, and applies a thin
# grid border (plus the existing header fill) to the used range on all
# three sheets.

$wb = $excel.ActiveWorkbook

$wsIndia = $wb.Worksheets.Item("India_LoginDatas")

# Grab the already-filled header format (blue fill) from India_LoginDatas
# before we touch anything, so it can be stamped onto the brand new sheets.
$wsIndia.Range("A1:D1").Copy() | Out-Null

function Fill-LoginSheet {
    param($ws, $country, $b3mail, $headerSource)

    # Header row
    $ws.Range("A1").Value = "country"
    $ws.Range("B1").Value = "emailID"
    $ws.Range("C1").Value = "Password"
    $ws.Range("D1").Value = "expected"

    if ($headerSource -ne $null) {
        $headerSource.Range("A1:D1").Copy() | Out-Null
        $ws.Range("A1:D1").PasteSpecial(-4122) | Out-Null
    }

    # Body rows
    $ws.Range("A2").Value = $country
    $ws.Range("B2").Value = "jkjkj@gmail.com"
    $ws.Range("C2").Value = "tes34"
    $ws.Range("D2").Value = "Invalid"

    $ws.Range("A3").Value = $country
    $ws.Range("B3").Value = $b3mail
    $ws.Range("C3").Value = "test1234"
    $ws.Range("D3").Value = "Valid"

    $ws.Range("A4").Value = $country
    $ws.Range("B4").Value = "mobile@gmail.com"
    $ws.Range("C4").Value = "test123"
    $ws.Range("D4").Value = "Invalid"

    $ws.Range("A5").Value = $country
    $ws.Range("B5").Value = "mobile1@gmail.com"
    $ws.Range("C5").Value = "test1234"
    $ws.Range("D5").Value = "Invalid"

    # Thin box border around the whole used range
    $ws.Range("A1:D5").Borders.LineStyle = 1
}

function AutoFit-Columns {
    param($ws)
    $ws.Columns.Item(1).AutoFit() | Out-Null
    $ws.Columns.Item(2).AutoFit() | Out-Null
    $ws.Columns.Item(3).AutoFit() | Out-Null
    $ws.Columns.Item(4).AutoFit() | Out-Null
}

function Set-LoginHyperlinks {
    param($ws, $b3mail, $b2mail, $b5mail, $b4mail)

    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("B3"), "mailto:$b3mail") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B2"), "mailto:$b2mail") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B5"), "mailto:$b5mail") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B4"), "mailto:$b4mail") | Out-Null
}

# ---------------------------------------------------------------------------
# India_LoginDatas (existing sheet) - refresh its data/format in place
# ---------------------------------------------------------------------------
Fill-LoginSheet $wsIndia "IND" "mobile@gmail.com" $null
Set-LoginHyperlinks $wsIndia "mobile@gmail.com" "jkjkj@gmail.com" "mobile1@gmail.com" "mobile@gmail.com"

$wsIndia.Range("D2").Select() | Out-Null

# ---------------------------------------------------------------------------
# UK_LoginDatas (new sheet, placed after India_LoginDatas)
# ---------------------------------------------------------------------------
$wsUK = $wb.Worksheets.Add($null, $wsIndia)
$wsUK.Name = "UK_LoginDatas"

Fill-LoginSheet $wsUK "UK" "priya@gmail.com" $wsIndia
Set-LoginHyperlinks $wsUK "priya@gmail.com" "jkjkj@gmail.com" "mobile1@gmail.com" "mobile@gmail.com"
AutoFit-Columns $wsUK

$wsUK.Range("B9").Select() | Out-Null
$wsIndia.Select() | Out-Null

# ---------------------------------------------------------------------------
# USA_LoginDatas (new sheet, placed after UK_LoginDatas)
# ---------------------------------------------------------------------------
$wsUSA = $wb.Worksheets.Add($null, $wsUK)
$wsUSA.Name = "USA_LoginDatas"

Fill-LoginSheet $wsUSA "USA" "build28@gmail.com" $wsIndia
Set-LoginHyperlinks $wsUSA "build28@gmail.com" "jkjkj@gmail.com" "mobile1@gmail.com" "mobile@gmail.com"
AutoFit-Columns $wsUSA

$wsUSA.Range("A11").Select() | Out-Null
$wsIndia.Select() | Out-Null
